$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (column D) and Volume(1h) (column E) values for each
# coin row with the latest figures from the source feed.
#
# Prices are stored as plain text (not numbers) throughout this sheet, since
# some of them contain multiple "." separators (e.g. "36.476.42"). For the
# rows whose new price happens to look like an ordinary decimal number, a
# leading apostrophe is used so Excel keeps storing the text as-is instead of
# silently converting it to a numeric cell; the cell style is then reset back
# to "Normal" so no extra text-formatting flag is left behind on the cell.

$ws.Range("D2").Value = "36.476.42"
$ws.Range("E2").Value = "  +1.22%  "

$ws.Range("D3").Value = "1.944.89"
$ws.Range("E3").Value = "  -0.90%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'243.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").Value = "'0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.63%  "

$ws.Range("D7").Value = "'58.28"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.96%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "'0.365"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.98%  "

$ws.Range("D10").Value = "'55.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "

$ws.Range("D11").Value = "'0.0836"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.70%  "

$ws.Range("D12").Value = "'0.104"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.58%  "

$ws.Range("D13").Value = "'21.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.93%  "

$ws.Range("D14").Value = "'0.819"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.09%  "

$ws.Range("D15").Value = "2.226.76"

$ws.Range("D16").Value = "'13.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.28%  "

$ws.Range("D17").Value = "'5.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.24%  "

$ws.Range("D18").Value = "1.953.77"
$ws.Range("E18").Value = "  -0.63%  "

$ws.Range("D19").Value = "36.339.85"
$ws.Range("E19").Value = "  +0.96%  "

$ws.Range("D20").Value = "'69.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.95%  "

$ws.Range("D21").Value = "0.0₃0861"
$ws.Range("E21").Value = "  +0.74%  "

$ws.Range("D22").Value = "'229.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.98%  "

$ws.Range("D23").Value = "'5.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.68%  "

$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").Value = "'2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.46%  "

$ws.Range("D26").Value = "'2.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").Value = "'9.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.56%  "

$ws.Range("D28").Value = "'161.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.90%  "

$ws.Range("D29").Value = "'19.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.08%  "

$ws.Range("D30").Value = "'0.127"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.35%  "

$ws.Range("E31").Value = "  -1.35%  "

$ws.Range("E32").Value = "  +1.10%  "

$ws.Range("D33").Value = "'4.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.40%  "

$ws.Range("D34").Value = "'0.0625"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.92%  "

$ws.Range("D35").Value = "'4.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.64%  "

$ws.Range("E36").Value = "  -1.80%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("E38").Value = "  -2.91%  "

$ws.Range("E39").Value = "  -6.24%  "

$ws.Range("D40").Value = "'3.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.22%  "

$ws.Range("D41").Value = "'0.0982"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("E43").Value = "  -4.26%  "

$ws.Range("E44").Value = "  -1.37%  "

$ws.Range("D45").Value = "'16.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").Value = "1.347.83"
$ws.Range("E46").Value = "  +0.77%  "

$ws.Range("E47").Value = "  -4.85%  "

$ws.Range("D48").Value = "'87.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.94%  "

$ws.Range("D49").Value = "'7.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.65%  "

$ws.Range("D50").Value = "'2.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.03%  "

$ws.Range("D51").Value = "'45.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.53%  "

